$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Weapons"
$ws.Range("B7").Value = "Movie"
$ws.Range("C7").Value = "Horror"
$ws.Range("D7").Value = "Medium"
$ws.Range("E7").Value = "Watched"
$ws.Range("F7").Value = 8.5
